# Updated symbol list on Fri Jan 20 20:57:11 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped values. Values are stored as literal text (leading
# apostrophe forces text so things like "297.29" and "1.04%" are kept
# exactly as scraped instead of becoming numeric/percentage cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.29"
$ws.Range("E2").Value = "'1.04%"
$ws.Range("D3").Value = "'31.63"
$ws.Range("E3").Value = "'2.02%"
$ws.Range("D4").Value = "'4.984"
$ws.Range("E4").Value = "'0.71%"
$ws.Range("D5").Value = "'0.07686"
$ws.Range("E5").Value = "'4.69%"
$ws.Range("D6").Value = "'2.237"
$ws.Range("E6").Value = "'-3.19%"
$ws.Range("D7").Value = "'7.896"
$ws.Range("E7").Value = "'1.95%"
$ws.Range("D8").Value = "'0.9243"
$ws.Range("E8").Value = "'1.77%"
$ws.Range("D9").Value = "'0.09916"
$ws.Range("E9").Value = "'22.90%"
$ws.Range("D10").Value = "'0.1747"
$ws.Range("E10").Value = "'4.03%"
$ws.Range("D11").Value = "'0.08389"
$ws.Range("E11").Value = "'2.71%"
$ws.Range("D12").Value = "'0.03240"
$ws.Range("E12").Value = "'4.57%"
$ws.Range("D13").Value = "'0.09826"
$ws.Range("E13").Value = "'-2.40%"
$ws.Range("D14").Value = "'0.001469"
$ws.Range("E14").Value = "'-3.30%"
$ws.Range("D15").Value = "'0.005736"
$ws.Range("E15").Value = "'0.18%"
$ws.Range("E16").Value = "'1.31%"
$ws.Range("D17").Value = "'3.787"
$ws.Range("E18").Value = "'5.82%"
$ws.Range("D19").Value = "'0.3365"
$ws.Range("E19").Value = "'1.15%"
$ws.Range("D20").Value = "'0.1326"
$ws.Range("E20").Value = "'1.72%"
$ws.Range("D21").Value = "'4.071"
$ws.Range("E21").Value = "'2.59%"
$ws.Range("D22").Value = "'0.2271"
$ws.Range("E22").Value = "'8.22%"
$ws.Range("D23").Value = "'0.04514"
$ws.Range("E23").Value = "'-0.79%"
$ws.Range("D24").Value = "'0.001210"
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("D25").Value = "'0.004365"
$ws.Range("E25").Value = "'-6.24%"
$ws.Range("D26").Value = "'0.0001285"
$ws.Range("E26").Value = "'-1.16%"
$ws.Range("D27").Value = "'0.0003356"
$ws.Range("E27").Value = "'-1.04%"
$ws.Range("D39").Value = "'0.01700"
$ws.Range("E39").Value = "'5.52%"
$ws.Range("D40").Value = "'0.04624"
$ws.Range("E40").Value = "'4.49%"
$ws.Range("D41").Value = "'0.007485"
$ws.Range("E41").Value = "'1.17%"
$ws.Range("D42").Value = "'0.009722"
$ws.Range("E42").Value = "'13.07%"
$ws.Range("D43").Value = "'0.1385"
$ws.Range("E43").Value = "'4.03%"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("E44").Value = "'5.96%"
$ws.Range("D45").Value = "'0.009584"
$ws.Range("E45").Value = "'0.81%"
$ws.Range("D46").Value = "'0.00006034"
$ws.Range("E46").Value = "'6.72%"
$ws.Range("D47").Value = "'0.00000000742"
$ws.Range("E47").Value = "'-1.03%"
$ws.Range("D48").Value = "'2.654"
$ws.Range("E48").Value = "'18.45%"
$ws.Range("D49").Value = "'0.001978"
$ws.Range("E49").Value = "'-31.71%"
$ws.Range("D50").Value = "'0.00002077"
$ws.Range("E50").Value = "'-1.03%"
$ws.Range("D51").Value = "'0.0001978"
$ws.Range("E51").Value = "'-1.03%"
